# Apply the "List 1" paragraph style to the two paragraphs that were
# touched in the source revision:
#   1. The paragraph beginning "    (b)(iii) Submit the required
#      contracting officer determination for ..."
#   2. The paragraph that starts right after the
#      "5212.302 Tailoring of provisions and clauses ..." heading and
#      begins with four spaces followed by "(c)  Tailoring inconsistent
#      with customary commercial practice."

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text

    if ($text.Contains("(b)(iii) Submit the required contracting officer determination for")) {
        $p.Range.set_Style("List 1")
    }
    elseif ($text.Contains("(c)  Tailoring inconsistent with customary commercial")) {
        $p.Range.set_Style("List 1")
    }
}
